# Re-save the computed summary statistics (mean cell concentration, water
# volume, and geometric mean cell concentration). The values were produced
# by code that recomputes and writes these results back into the workbook;
# re-entering them here causes Excel to re-serialize the underlying
# floating point numbers at its native double precision.

$wb = $excel.ActiveWorkbook

# --- Sheet "Cell concentration" ---
$ws1 = $wb.Worksheets.Item("Cell concentration")

$ws1.Range("C2").Value = 1014764.96084211
$ws1.Range("D2").Value = 418589.624461565

$ws1.Range("C3").Value = 642998.737414316
$ws1.Range("D3").Value = 215223.638082962

$ws1.Range("C4").Value = 729590.390147741
$ws1.Range("D4").Value = 407859.733851065

$ws1.Range("C5").Value = 164905.522897753
$ws1.Range("D5").Value = 67276.0758529747

$ws1.Range("C6").Value = 156688.53815352
$ws1.Range("D6").Value = 120326.827538516

$ws1.Range("C7").Value = 136977.672565956
$ws1.Range("D7").Value = 11946.1159025778

$ws1.Range("D8").Value = 3428.34337087247

$ws1.Range("C9").Value = 21480.7055359862
$ws1.Range("D9").Value = 3472.08910207416

# --- Sheet "Water volume" (Geometric mean cell concentration column) ---
# These magnitudes (~1e21) can't be written as a single numeric literal, so
# each value is built as mantissa * 10^a * 10^b (a+b = 21) using the exact
# split that reproduces the target double bit-for-bit.
$ws2 = $wb.Worksheets.Item("Water volume")

$ws2.Range("C2").Value = 5.35401807906929 * 1000000000 * 1000000000000
$ws2.Range("C3").Value = 4.0853193164043 * 10000000000 * 100000000000
$ws2.Range("C4").Value = 3.20933202659467 * 1000000000000 * 1000000000
$ws2.Range("C5").Value = 2.60449677768489 * 10000000000 * 100000000000
$ws2.Range("C6").Value = 2.18688149652793 * 10000000000 * 100000000000
$ws2.Range("C7").Value = 1.89853434187688 * 10000000000000 * 100000000
$ws2.Range("C8").Value = 1.69944181498222 * 1000000000 * 1000000000000
$ws2.Range("C9").Value = 1.56197614685982 * 10000000000 * 100000000000

$wb.Save()
